# 8.10.1 — add 2023 column (T), refresh population-derived rates for
# 2010-2020 (columns G:Q) after the underlying adult-population figures
# (row 8) were revised, and drop the now-stale shared formulas in rows 4-5
# in favour of static values (matching columns R:S, which were already
# literals).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 4 (branches per 100 000 adults): replace formulas D4:Q4 with
#    static values. D4:F4 keep their old figures (G8 denominator for
#    those years is unchanged); G4:Q4 get the values recomputed against
#    the new row-8 population figures.
# ---------------------------------------------------------------------
$row4Cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
$row4Vals = @(
    6.1074687240787666,
    6.5454292116044552,
    6.6165964726065987,
    5.6686326881838296,
    6.7294661864194607,
    7.2207098269445202,
    7.3191488059459031,
    7.3364889416826751,
    7.7198339498137045,
    7.8258279858854918,
    7.5789073543911334,
    7.4985248229203512,
    7.4141082446031374,
    7.0384645318913508
)
for ($i = 0; $i -lt $row4Cols.Length; $i++) {
    $ws.Range("$($row4Cols[$i])4").Value = $row4Vals[$i]
}

# ---------------------------------------------------------------------
# 2. Row 5 (ATMs per 100 000 adults): same treatment.
# ---------------------------------------------------------------------
$row5Cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
$row5Vals = @(
    2.4713943209062914,
    4.8880029305544008,
    8.7206187818873584,
    7.5223561738932325,
    12.187222227373827,
    15.844003577108481,
    20.86084920997822,
    24.989132982250201,
    30.387626630476873,
    31.39992710386154,
    33.570520663807748,
    36.978706525491944,
    39.297069085946042,
    41.869840292276756
)
for ($i = 0; $i -lt $row5Cols.Length; $i++) {
    $ws.Range("$($row5Cols[$i])5").Value = $row5Vals[$i]
}

# ---------------------------------------------------------------------
# 3. Row 8 (adult population): revised figures for 2010-2020 (G:Q).
# ---------------------------------------------------------------------
$row8Cols = @("G","H","I","J","K","L","M","N","O","P","Q")
$row8Vals = @(3722238, 3774445, 3850037, 3921221, 3993736, 4067445, 4140137, 4209050, 4280842, 4356559, 4432785)
for ($i = 0; $i -lt $row8Cols.Length; $i++) {
    $ws.Range("$($row8Cols[$i])8").Value = $row8Vals[$i]
}

# ---------------------------------------------------------------------
# 4. New column T (2023). Copy the formatting from column S (the prior
#    last column) across the data rows, then fill in the 2023 figures.
# ---------------------------------------------------------------------
$ws.Range("S3:S8").Copy() | Out-Null
$ws.Range("T3:T8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("T3").Value = 2023
$ws.Range("T4").Value = 6.7904451646088795
$ws.Range("T5").Value = 47.957518975050206
$ws.Range("T6").Value = 320
$ws.Range("T7").Value = 2260
$ws.Range("T8").Value = 4712504

# ---------------------------------------------------------------------
# 5. Row-height tweaks that came along with the new column.
# ---------------------------------------------------------------------
$ws.Rows(4).RowHeight = 27
$ws.Rows(6).RowHeight = 15.75
